{"js": "// Apply the Gutachten.Sachverhalt2 text corrections.\n//\n// Strategy: locate each target text run with body.search() (exact, case-sensitive;\n// every search string below is the full content of a single <w:t> run, so it is\n// unique in the document) and rewrite it with\n// insertText(..., Word.InsertLocation.replace). A vertical-tab character (\"\\v\")\n// inside a replacement string becomes a manual line break (<w:br/>), matching how\n// Word represents <w:br/> in its text() - so a run that must be split across one or\n// two new line breaks is simply given a replacement string containing \"\\v\" (or\n// \"\\v\\v\" for two consecutive breaks).\n//\n// Two of the edits reorder/merge adjacent sentences; for those we locate both\n// sentences and build a single range spanning from the start of the first to the\n// end of the second with Range.expandTo(), then replace the whole span in one shot\n// so the sentence order (and break count) comes out exactly right.\n//\n// Important: each search()+insertText() is followed by its own context.sync().\n// The whole document body lives in a single <w:r>, so earlier text-length changes\n// shift the character offsets backing any Range object that was resolved before\n// that edit; syncing immediately after every mutation (and re-searching afterwards\n// for the next edit) keeps every subsequent range accurate.\n\nasync function replaceOnce(oldText, newText) {\n  const r = context.document.body.search(oldText, { matchCase: true });\n  r.load(\"text\");\n  await context.sync();\n  r.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(\n  \"Die Anordnung des Landratsamts k\u00f6nnte auf \u00a71 Abs.\",\n  \"Die Anordnung des Landratsamts k\u00f6nnte auf \u00a7 1 Abs.\"\n);\n\nawait replaceOnce(\n  \"1 in Verbindung mit \u00a7 7 Abs.1 Satz 1 DSchG gest\u00fctzt werden.\",\n  \"1 in Verbindung mit \u00a7 7 Abs.\\v\\v1 Satz 1 DSchG gest\u00fctzt werden.\"\n);\n\nawait replaceOnce(\n  \"Das Fachwerkhaus m\u00fcsste ein Kulturdenkmal sein.\",\n  \"Das Fachwerkhaus muss ein Kulturdenkmal sein.\"\n);\n\nawait replaceOnce(\n  \"Es k\u00f6nnte eine Gef\u00e4hrdung beim Kulturdenkmal (Fachwerkhaus) vorliegen.\",\n  \"Es k\u00f6nnte eine Gef\u00e4hrdung des Kulturdenkmals (Fachwerkhaus) vorliegen.\"\n);\n\n// \"Der G.K ist ebenfalls ...\" / \"Letztlich ist F.K ...\" swap: the old second\n// sentence is dropped and a brand-new sentence about F.K takes its place in front\n// of the (retained) G.K sentence.\n{\n  const a = context.document.body.search(\n    \"Der G.K ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.\",\n    { matchCase: true }\n  );\n  a.load(\"text\");\n  const b = context.document.body.search(\n    \"Letztlich ist F.K der richtige Pflichtige, da er \u00fcber ein h\u00f6heres Einkommen verf\u00fcgt und somit leistungsf\u00e4higer ist.\",\n    { matchCase: true }\n  );\n  b.load(\"text\");\n  await context.sync();\n\n  const span = a.items[0].expandTo(b.items[0]);\n  span.insertText(\n    \"Der F.K ist Eigent\u00fcmer des Fachwerkhauses, von dessen Dach eine Gef\u00e4hrdung f\u00fcr das Denkmal ausgeht.\\v\\vDer G.K ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\nawait replaceOnce(\n  \"Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen.\",\n  \"Die Anordnung des Landratsamts ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen.\"\n);\n\nawait replaceOnce(\n  \"Hier k\u00f6nnte \u00a72038 Abs.1 Satz 1 BGB relevant sein, da G.K als Miterbe allen Ma\u00dfnahmen zustimmen m\u00fcsste, dies jedoch verweigert.\",\n  \"Hier k\u00f6nnte \u00a7 2038 Abs.\\v\\v1 Satz 1 BGB relevant sein, da G.K als Miterbe allen Ma\u00dfnahmen zustimmen m\u00fcsste, dies jedoch verweigert.\"\n);\n\nawait replaceOnce(\n  \"Eine Ausnahme besteht gem\u00e4\u00df \u00a72038 Abs.1 Satz 2 Halbsatz 2 BGB, wenn die BSD als notwendige Erhaltungsma\u00dfnahme anzusehen ist.\",\n  \"Eine Ausnahme besteht gem\u00e4\u00df \u00a7 2038 Abs.\\v\\v1 Satz 2 Halbsatz 2 BGB, wenn die Reparatur als notwendige Erhaltungsma\u00dfnahme anzusehen ist.\"\n);\n\nawait replaceOnce(\"Formelle Vorrausetzung\", \"Formelle Voraussetzung\");\n\nawait replaceOnce(\"3 DSchG und \u00a7 46 Abs.\", \"3 DSchG und 46 Abs.\");\n\nawait replaceOnce(\n  \"\u00d6rtlich zust\u00e4ndig ist das Landratsamt Ortenaukreis gem\u00e4\u00df \u00a7 3 Abs.1 Nr.1 LVwVfG.\",\n  \"\u00d6rtlich zust\u00e4ndig ist das Landratsamt Ortenaukreis gem\u00e4\u00df \u00a7 3 Abs.\\v\\v1 Nr.\\v\\v1 LVwVfG.\"\n);\n\n// \"G.K ist beteiligt ...\" / \"(\u00a7903 BGB) Es besteht ...\" merge: the \"(\u00a7903 BGB)\"\n// parenthetical moves to the end of the first sentence.\n{\n  const a = context.document.body.search(\n    \"G.K ist beteiligt, da er Eigent\u00fcmer des Fachwerkhauses ist.\",\n    { matchCase: true }\n  );\n  a.load(\"text\");\n  const b = context.document.body.search(\n    \"(\u00a7903 BGB) Es besteht also ein rechtliches Interesse nach \u00a7 13 Abs.\",\n    { matchCase: true }\n  );\n  b.load(\"text\");\n  await context.sync();\n\n  const span = a.items[0].expandTo(b.items[0]);\n  span.insertText(\n    \"G.K ist beteiligt, da er Eigent\u00fcmer des Fachwerkhauses ist (\u00a7 903 BGB).\\v\\vEs besteht also ein rechtliches Interesse nach \u00a7 13 Abs.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\nawait replaceOnce(\n  \"Es besteht keine Befangenheit, da der Mitarbeiter, der den Fall bearbeitet, nach Sachverhalt nicht bekannt ist.\",\n  \"Es k\u00f6nnte eine Problematik im Sinne des \u00a7 21 LVwVfG gegeben sein, da F.K in der Vergangenheit das Landratsamt in Leserbriefen angegriffen hat.\\v\\vEine Befangenheit kann jedoch nicht festgestellt werden, da der Mitarbeiter, der den Fall bearbeitet, den Sachverhalt nicht kennt.\"\n);\n\nawait replaceOnce(\"Nach \u00a73 Abs.\", \"Nach \u00a7 3 Abs.\");\n\nawait replaceOnce(\"Rechtbehelfsbelehrung\", \"Rechtsbehelfsbelehrung\");\n", "ps1": "# Apply the Gutachten.Sachverhalt2 text corrections via Word COM interop.\n#\n# Strategy: for a plain 1-for-1 text rewrite, run Find/Replace on $d.Content with\n# Find.Text / Find.Replacement.Text (every search string below is the full content\n# of a single <w:t> run, so it is unambiguous). A backtick-v (\"`v\", vertical tab,\n# 0x0B) inside the replacement text becomes a manual line break (<w:br/>) - the same\n# character Word itself uses for a line break - so a run that must be split across\n# one or two new line breaks is simply given a replacement string containing \"`v\"\n# (or \"`v`v\" for two consecutive breaks).\n#\n# Two of the edits reorder/merge adjacent sentences; for those we find both\n# sentences to get the Start of the first and the End of the second, build a single\n# Range spanning both with $d.Range(start, end), and set that Range's .Text in one\n# shot so the sentence order (and break count) comes out exactly right.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n}\n\nReplace-Once \"Die Anordnung des Landratsamts k\u00f6nnte auf \u00a71 Abs.\" \"Die Anordnung des Landratsamts k\u00f6nnte auf \u00a7 1 Abs.\"\n\nReplace-Once \"1 in Verbindung mit \u00a7 7 Abs.1 Satz 1 DSchG gest\u00fctzt werden.\" \"1 in Verbindung mit \u00a7 7 Abs.`v`v1 Satz 1 DSchG gest\u00fctzt werden.\"\n\nReplace-Once \"Das Fachwerkhaus m\u00fcsste ein Kulturdenkmal sein.\" \"Das Fachwerkhaus muss ein Kulturdenkmal sein.\"\n\nReplace-Once \"Es k\u00f6nnte eine Gef\u00e4hrdung beim Kulturdenkmal (Fachwerkhaus) vorliegen.\" \"Es k\u00f6nnte eine Gef\u00e4hrdung des Kulturdenkmals (Fachwerkhaus) vorliegen.\"\n\n# \"Der G.K ist ebenfalls ...\" / \"Letztlich ist F.K ...\" swap: the old second\n# sentence is dropped and a brand-new sentence about F.K takes its place in front of\n# the (retained) G.K sentence.\n$rA = $d.Content\n$rA.Find.Execute(\"Der G.K ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.\") | Out-Null\n$startPos = $rA.Start\n$rB = $d.Content\n$rB.Find.Execute(\"Letztlich ist F.K der richtige Pflichtige, da er \u00fcber ein h\u00f6heres Einkommen verf\u00fcgt und somit leistungsf\u00e4higer ist.\") | Out-Null\n$endPos = $rB.End\n$span = $d.Range($startPos, $endPos)\n$span.Text = \"Der F.K ist Eigent\u00fcmer des Fachwerkhauses, von dessen Dach eine Gef\u00e4hrdung f\u00fcr das Denkmal ausgeht.`v`vDer G.K ist ebenfalls Eigent\u00fcmer des Fachwerkhauses und somit nach denselben Vorschriften pflichtig.\"\n\nReplace-Once \"Die Anordnung der BSD ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen.\" \"Die Anordnung des Landratsamts ist verh\u00e4ltnism\u00e4\u00dfig und ermessensgerecht, da eine kosteng\u00fcnstigere Reparatur mit Eternitplatten nicht geeignet ist, die Denkmalanforderungen zu erf\u00fcllen.\"\n\nReplace-Once \"Hier k\u00f6nnte \u00a72038 Abs.1 Satz 1 BGB relevant sein, da G.K als Miterbe allen Ma\u00dfnahmen zustimmen m\u00fcsste, dies jedoch verweigert.\" \"Hier k\u00f6nnte \u00a7 2038 Abs.`v`v1 Satz 1 BGB relevant sein, da G.K als Miterbe allen Ma\u00dfnahmen zustimmen m\u00fcsste, dies jedoch verweigert.\"\n\nReplace-Once \"Eine Ausnahme besteht gem\u00e4\u00df \u00a72038 Abs.1 Satz 2 Halbsatz 2 BGB, wenn die BSD als notwendige Erhaltungsma\u00dfnahme anzusehen ist.\" \"Eine Ausnahme besteht gem\u00e4\u00df \u00a7 2038 Abs.`v`v1 Satz 2 Halbsatz 2 BGB, wenn die Reparatur als notwendige Erhaltungsma\u00dfnahme anzusehen ist.\"\n\nReplace-Once \"Formelle Vorrausetzung\" \"Formelle Voraussetzung\"\n\nReplace-Once \"3 DSchG und \u00a7 46 Abs.\" \"3 DSchG und 46 Abs.\"\n\nReplace-Once \"\u00d6rtlich zust\u00e4ndig ist das Landratsamt Ortenaukreis gem\u00e4\u00df \u00a7 3 Abs.1 Nr.1 LVwVfG.\" \"\u00d6rtlich zust\u00e4ndig ist das Landratsamt Ortenaukreis gem\u00e4\u00df \u00a7 3 Abs.`v`v1 Nr.`v`v1 LVwVfG.\"\n\n# \"G.K ist beteiligt ...\" / \"(\u00a7903 BGB) Es besteht ...\" merge: the \"(\u00a7903 BGB)\"\n# parenthetical moves to the end of the first sentence.\n$rC = $d.Content\n$rC.Find.Execute(\"G.K ist beteiligt, da er Eigent\u00fcmer des Fachwerkhauses ist.\") | Out-Null\n$startPos2 = $rC.Start\n$rD = $d.Content\n$rD.Find.Execute(\"(\u00a7903 BGB) Es besteht also ein rechtliches Interesse nach \u00a7 13 Abs.\") | Out-Null\n$endPos2 = $rD.End\n$span2 = $d.Range($startPos2, $endPos2)\n$span2.Text = \"G.K ist beteiligt, da er Eigent\u00fcmer des Fachwerkhauses ist (\u00a7 903 BGB).`v`vEs besteht also ein rechtliches Interesse nach \u00a7 13 Abs.\"\n\nReplace-Once \"Es besteht keine Befangenheit, da der Mitarbeiter, der den Fall bearbeitet, nach Sachverhalt nicht bekannt ist.\" \"Es k\u00f6nnte eine Problematik im Sinne des \u00a7 21 LVwVfG gegeben sein, da F.K in der Vergangenheit das Landratsamt in Leserbriefen angegriffen hat.`v`vEine Befangenheit kann jedoch nicht festgestellt werden, da der Mitarbeiter, der den Fall bearbeitet, den Sachverhalt nicht kennt.\"\n\nReplace-Once \"Nach \u00a73 Abs.\" \"Nach \u00a7 3 Abs.\"\n\nReplace-Once \"Rechtbehelfsbelehrung\" \"Rechtsbehelfsbelehrung\"\n"}
